# TC_134 test data sheet update:
# On sheet "Add Devices Loop A", column I (Isolator Units) rows 8-10 were
# re-labelled from the old placeholder text "801 H - 1" to "Isolator Units".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Add Devices Loop A")

$ws1.Range("I8").Value  = "Isolator Units"
$ws1.Range("I9").Value  = "Isolator Units"
$ws1.Range("I10").Value = "Isolator Units"
